$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Prok1"
$ws.Range("C2").Value = "Prokr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1903723333333333
$ws.Range("H2").Value = 0.571117
$ws.Range("I2").Value = 0.7691179579859432
$ws.Range("J2").Value = 0.7691179579859432
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.145415666666666
$ws.Range("N2").Value = 18.436247
$ws.Range("O2").Value = 0.9968753168251276
$ws.Range("P2").Value = 0.9968753168251278
$ws.Range("Q2").Value = 1.169917119766555
$ws.Range("R2").Value = 10.529254077899
$ws.Range("S2").Value = 0.7667147080431324
$ws.Range("T2").Value = 0.7667147080431324

# --- New row 3 ---
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Prok1"
$ws.Range("C3").Value = "Prokr2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1903723333333333
$ws.Range("H3").Value = 0.571117
$ws.Range("I3").Value = 0.7691179579859432
$ws.Range("J3").Value = 0.7691179579859432
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01926266666666667
$ws.Range("N3").Value = 0.057788
$ws.Range("O3").Value = 0.003124683174872331
$ws.Range("P3").Value = 0.003124683174872331
$ws.Range("Q3").Value = 0.003667078799555556
$ws.Range("R3").Value = 0.033003709196
$ws.Range("S3").Value = 0.002403249942810841
$ws.Range("T3").Value = 0.002403249942810841

# --- New row 4 ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Prok1"
$ws.Range("C4").Value = "Prokr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.057148
$ws.Range("H4").Value = 0.171444
$ws.Range("I4").Value = 0.2308820420140567
$ws.Range("J4").Value = 0.2308820420140567
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.145415666666666
$ws.Range("N4").Value = 18.436247
$ws.Range("O4").Value = 0.9968753168251276
$ws.Range("P4").Value = 0.9968753168251278
$ws.Range("Q4").Value = 0.3511982145186666
$ws.Range("R4").Value = 3.160783930667999
$ws.Range("S4").Value = 0.2301606087819952
$ws.Range("T4").Value = 0.2301606087819953

# --- New row 5 ---
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Prok1"
$ws.Range("C5").Value = "Prokr2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.057148
$ws.Range("H5").Value = 0.171444
$ws.Range("I5").Value = 0.2308820420140567
$ws.Range("J5").Value = 0.2308820420140567
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01926266666666667
$ws.Range("N5").Value = 0.057788
$ws.Range("O5").Value = 0.003124683174872331
$ws.Range("P5").Value = 0.003124683174872331
$ws.Range("Q5").Value = 0.001100822874666667
$ws.Range("R5").Value = 0.009907405871999999
$ws.Range("S5").Value = 0.0007214332320614897
$ws.Range("T5").Value = 0.0007214332320614898
